# Apply updated crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.084.68"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "1.638.06"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.82"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3945"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3878"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.30"
$ws.Range("E9").Value = "  +4.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.389"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9981"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08530"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.155"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001306"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.689"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "1.621.71"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.94"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06928"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.22"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.901"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9970"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.54"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("D24").Value = "24.075.44"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.452"
$ws.Range("E25").Value = "  +4.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.894"
$ws.Range("E26").Value = "  +3.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.41"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.47"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "141.38"
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.380"
$ws.Range("E30").Value = "  -5.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.016"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.541"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").Value = "1.808.80"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.014"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08194"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02936"
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.704"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2705"
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09220"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.44"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.90"
$ws.Range("E41").Value = "  +3.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7643"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.429"
$ws.Range("E43").Value = "  -3.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.25"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6991"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.494"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.110"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9977"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08345"
$ws.Range("E49").Value = "  -3.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.62"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.422"
$ws.Range("E51").Value = "  +16.79%  "
